$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.165.53'
$ws.Range("E2").Value = '  -1.75%  '

$ws.Range("D3").Value = '2.246.97'
$ws.Range("E3").Value = '  -1.80%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.82'
$ws.Range("E5").Value = '  -1.78%  '

$ws.Range("E6").Value = '  -2.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '77.42'
$ws.Range("E7").Value = '  +4.67%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.634'
$ws.Range("E9").Value = '  -1.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.02'
$ws.Range("E10").Value = '  +7.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0958'
$ws.Range("E11").Value = '  -2.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.23'
$ws.Range("E12").Value = '  -2.71%  '

$ws.Range("E13").Value = '  -2.35%  '

$ws.Range("D14").Value = '2.582.89'
$ws.Range("E14").Value = '  -1.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.87'
$ws.Range("E15").Value = '  -2.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.863'
$ws.Range("E16").Value = '  -1.44%  '

$ws.Range("D17").Value = '2.223.58'

$ws.Range("D18").Value = '42.039.22'
$ws.Range("E18").Value = '  -1.77%  '

$ws.Range("D19").Value = '0.0₃0986'
$ws.Range("E19").Value = '  -1.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.15'
$ws.Range("E20").Value = '  -2.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.02'
$ws.Range("E21").Value = '  -0.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.32'
$ws.Range("E22").Value = '  +4.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.19'
$ws.Range("E23").Value = '  -2.25%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("E25").Value = '  -1.71%  '

$ws.Range("E26").Value = '  -6.31%  '

$ws.Range("E27").Value = '  -5.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.31'
$ws.Range("E28").Value = '  +12.82%  '

$ws.Range("E29").Value = '  -1.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.61'
$ws.Range("E30").Value = '  +1.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.59'
$ws.Range("E31").Value = '  -2.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.98'
$ws.Range("E32").Value = '  +9.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0833'
$ws.Range("E33").Value = '  +0.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.122'
$ws.Range("E34").Value = '  -4.67%  '

$ws.Range("E35").Value = '  -0.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.55'
$ws.Range("E36").Value = '  -1.66%  '

$ws.Range("E37").Value = '  +2.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.35'
$ws.Range("E38").Value = '  -0.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0303'
$ws.Range("E39").Value = '  -2.35%  '

$ws.Range("E40").Value = '  +0.22%  '

$ws.Range("E41").Value = '  -6.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '113.44'
$ws.Range("E42").Value = '  +12.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.204'
$ws.Range("E43").Value = '  -5.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.19'
$ws.Range("E44").Value = '  -1.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.71'
$ws.Range("E45").Value = '  -4.96%  '

$ws.Range("E46").Value = '  -2.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.998'
$ws.Range("E47").Value = '  -0.31%  '

$ws.Range("E48").Value = '  -2.68%  '

$ws.Range("E49").Value = '  -1.05%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("E50").Value = '  -0.62%  '

$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.22'
$ws.Range("E51").Value = '  -13.44%  '
